$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-330 currently holds the date serial
# 45190 (2023-09-21). Update it to 45192 (2023-09-23) for every data row.
$ws.Range("C2:C330").Value = 45192
